$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format + borders) of the last existing data
# row (row 8) down into the new row 9, matching how the workbook's
# formatting pattern was extended for the newly appended data.
$ws.Range("A8:M8").Copy() | Out-Null
$ws.Range("A9:M9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A9").Value = 45756
$ws.Range("B9").Value = 36.1
$ws.Range("C9").Value = 40.6
$ws.Range("D9").Value = 38.9
$ws.Range("E9").Value = 39.7
$ws.Range("F9").Value = 38.6
$ws.Range("G9").Value = 33
$ws.Range("H9").Value = 31.8
$ws.Range("I9").Value = 38.4
$ws.Range("J9").Value = 31.9
$ws.Range("K9").Value = 33.7
$ws.Range("L9").Value = 28.2
$ws.Range("M9").Value = 34.3

# Update the selected cell to match the new active cell in the worksheet view.
$ws.Range("I6").Select()
